# Applies the scraped-data refresh for the Lebanon Premier League 2023-2024 sheet:
#  1) Rows 16 and 17 had their match details (teams, odds, timestamps, url) swapped
#     back to the correct order (Al Sahel vs Al Ansar / Tripoli vs Al Ghazieh).
#  2) Two new matches were appended as rows 40 and 41 (Safa vs Tadamon, Tripoli vs Al Ahed).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Swap the content of rows 16 and 17 -------------------------------
# Columns G, I, K, O, S are identical between the two rows, so only the
# remaining columns need to be exchanged.

$row16 = @{
    "F" = "Tripoli"
    "H" = "Al Ghazieh"
    "J" = 2
    "L" = 1.67
    "M" = "20/08/2023 15:00"
    "N" = 3.08
    "P" = 3.64
    "Q" = "20/08/2023 15:00"
    "R" = 3.46
    "T" = 5.06
    "U" = "20/08/2023 15:00"
    "V" = "https://www.betexplorer.com/football/lebanon/premier-league/tripoli-sc-al-ghazieh/QDWlcvKR/"
}

$row17 = @{
    "F" = "Al Sahel"
    "H" = "Al Ansar"
    "J" = 3.98
    "L" = 4.23
    "M" = "20/08/2023 14:06"
    "N" = 3.39
    "P" = 3.49
    "Q" = "20/08/2023 15:31"
    "R" = 1.76
    "T" = 1.82
    "U" = "20/08/2023 14:06"
    "V" = "https://www.betexplorer.com/football/lebanon/premier-league/al-sahel-al-ansar/SMIGhbSr/"
}

foreach ($col in $row16.Keys) {
    $ws.Range("$col" + "16").Value = $row16[$col]
}

foreach ($col in $row17.Keys) {
    $ws.Range("$col" + "17").Value = $row17[$col]
}

# --- 2) Append the two new rows (40 and 41) ------------------------------
# Copy the formatting of the last existing row (39) onto the new rows first
# (this gives A40/A41 the bold+bordered "index" style and E40/E41 the
# datetime number format, matching the rest of the sheet).

$ws.Range("A39:V39").Copy()
$ws.Range("A40:V41").PasteSpecial(-4122)

$row40 = @{
    "A" = 39
    "B" = "lebanon"
    "C" = "premier-league"
    "D" = "2023-2024"
    "E" = 45234.55208333334
    "F" = "Safa"
    "G" = 4
    "H" = "Tadamon"
    "I" = 1
    "J" = 1.69
    "K" = "03/11/2023 01:43"
    "L" = 1.78
    "M" = "04/11/2023 10:21"
    "N" = 3.4
    "O" = "03/11/2023 01:43"
    "P" = 3.47
    "Q" = "04/11/2023 11:20"
    "R" = 4.32
    "S" = "03/11/2023 01:43"
    "T" = 4.42
    "U" = "04/11/2023 10:21"
    "V" = "https://www.betexplorer.com/football/lebanon/premier-league/safa-tadamon/8Yu0Lm2n/"
}

$row41 = @{
    "A" = 40
    "B" = "lebanon"
    "C" = "premier-league"
    "D" = "2023-2024"
    "E" = 45234.55208333334
    "F" = "Tripoli"
    "G" = 0
    "H" = "Al Ahed"
    "I" = 1
    "J" = 8.61
    "K" = "03/11/2023 01:43"
    "L" = 13.27
    "M" = "04/11/2023 13:13"
    "N" = 5.91
    "O" = "03/11/2023 01:43"
    "P" = 6.72
    "Q" = "04/11/2023 13:13"
    "R" = 1.21
    "S" = "03/11/2023 01:43"
    "T" = 1.18
    "U" = "04/11/2023 11:17"
    "V" = "https://www.betexplorer.com/football/lebanon/premier-league/tripoli-sc-al-ahed/UmU8JRWb/"
}

foreach ($col in $row40.Keys) {
    $ws.Range("$col" + "40").Value = $row40[$col]
}

foreach ($col in $row41.Keys) {
    $ws.Range("$col" + "41").Value = $row41[$col]
}

Write-Output "Edit applied."
